$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 67.91996
$ws.Cells.Item(2, 8).Value = 203.75988
$ws.Cells.Item(2, 9).Value = 0.02375577759132129
$ws.Cells.Item(2, 10).Value = 0.02375577759132129
$ws.Cells.Item(2, 13).Value = 6.066605666666667
$ws.Cells.Item(2, 14).Value = 18.199817
$ws.Cells.Item(2, 15).Value = 0.849784628791665
$ws.Cells.Item(2, 16).Value = 0.8497846287916652
$ws.Cells.Item(2, 17).Value = 412.0436142157733
$ws.Cells.Item(2, 18).Value = 3708.39252794196
$ws.Cells.Item(2, 19).Value = 0.02018729464209832
$ws.Cells.Item(2, 20).Value = 0.02018729464209832

$ws.Cells.Item(3, 7).Value = 67.91996
$ws.Cells.Item(3, 8).Value = 203.75988
$ws.Cells.Item(3, 9).Value = 0.02375577759132129
$ws.Cells.Item(3, 10).Value = 0.02375577759132129
$ws.Cells.Item(3, 15).Value = 0.1196497582104962
$ws.Cells.Item(3, 16).Value = 0.1196497582104962
$ws.Cells.Item(3, 17).Value = 58.01578087285333
$ws.Cells.Item(3, 18).Value = 522.14202785568
$ws.Cells.Item(3, 19).Value = 0.002842373044903915
$ws.Cells.Item(3, 20).Value = 0.002842373044903917

$ws.Cells.Item(4, 7).Value = 67.91996
$ws.Cells.Item(4, 8).Value = 203.75988
$ws.Cells.Item(4, 9).Value = 0.02375577759132129
$ws.Cells.Item(4, 10).Value = 0.02375577759132129
$ws.Cells.Item(4, 13).Value = 0.1824346666666667
$ws.Cells.Item(4, 14).Value = 0.547304
$ws.Cells.Item(4, 15).Value = 0.02555468148257719
$ws.Cells.Item(4, 16).Value = 0.02555468148257719
$ws.Cells.Item(4, 17).Value = 12.39095526261333
$ws.Cells.Item(4, 18).Value = 111.51859736352
$ws.Cells.Item(4, 19).Value = 0.0006070713297171602
$ws.Cells.Item(4, 20).Value = 0.0006070713297171604

$ws.Cells.Item(5, 7).Value = 67.91996
$ws.Cells.Item(5, 8).Value = 203.75988
$ws.Cells.Item(5, 9).Value = 0.02375577759132129
$ws.Cells.Item(5, 10).Value = 0.02375577759132129
$ws.Cells.Item(5, 13).Value = 0.035773
$ws.Cells.Item(5, 14).Value = 0.107319
$ws.Cells.Item(5, 15).Value = 0.005010931515261538
$ws.Cells.Item(5, 16).Value = 0.005010931515261539
$ws.Cells.Item(5, 17).Value = 2.42970072908
$ws.Cells.Item(5, 18).Value = 21.86730656172
$ws.Cells.Item(5, 19).Value = 0.0001190385746018957
$ws.Cells.Item(5, 20).Value = 0.0001190385746018957

$ws.Cells.Item(6, 9).Value = 0.9176057312269553
$ws.Cells.Item(6, 10).Value = 0.9176057312269554
$ws.Cells.Item(6, 13).Value = 6.066605666666667
$ws.Cells.Item(6, 14).Value = 18.199817
$ws.Cells.Item(6, 15).Value = 0.849784628791665
$ws.Cells.Item(6, 16).Value = 0.8497846287916652
$ws.Cells.Item(6, 17).Value = 15915.85796198022
$ws.Cells.Item(6, 18).Value = 143242.7216578219
$ws.Cells.Item(6, 19).Value = 0.7797672456878026
$ws.Cells.Item(6, 20).Value = 0.7797672456878028

$ws.Cells.Item(7, 9).Value = 0.9176057312269553
$ws.Cells.Item(7, 10).Value = 0.9176057312269554
$ws.Cells.Item(7, 15).Value = 0.1196497582104962
$ws.Cells.Item(7, 16).Value = 0.1196497582104962
$ws.Cells.Item(7, 19).Value = 0.1097913038738707
$ws.Cells.Item(7, 20).Value = 0.1097913038738708

$ws.Cells.Item(8, 9).Value = 0.9176057312269553
$ws.Cells.Item(8, 10).Value = 0.9176057312269554
$ws.Cells.Item(8, 13).Value = 0.1824346666666667
$ws.Cells.Item(8, 14).Value = 0.547304
$ws.Cells.Item(8, 15).Value = 0.02555468148257719
$ws.Cells.Item(8, 16).Value = 0.02555468148257719
$ws.Cells.Item(8, 17).Value = 478.6208963542667
$ws.Cells.Item(8, 18).Value = 4307.588067188401
$ws.Cells.Item(8, 19).Value = 0.02344912218809217
$ws.Cells.Item(8, 20).Value = 0.02344912218809218

$ws.Cells.Item(9, 9).Value = 0.9176057312269553
$ws.Cells.Item(9, 10).Value = 0.9176057312269554
$ws.Cells.Item(9, 13).Value = 0.035773
$ws.Cells.Item(9, 14).Value = 0.107319
$ws.Cells.Item(9, 15).Value = 0.005010931515261538
$ws.Cells.Item(9, 16).Value = 0.005010931515261539
$ws.Cells.Item(9, 17).Value = 93.85116128485001
$ws.Cells.Item(9, 18).Value = 844.66045156365
$ws.Cells.Item(9, 19).Value = 0.004598059477189759
$ws.Cells.Item(9, 20).Value = 0.00459805947718976

$ws.Cells.Item(10, 7).Value = 1.376679
$ws.Cells.Item(10, 8).Value = 4.130037
$ws.Cells.Item(10, 9).Value = 0.0004815091195378001
$ws.Cells.Item(10, 10).Value = 0.0004815091195378002
$ws.Cells.Item(10, 13).Value = 6.066605666666667
$ws.Cells.Item(10, 14).Value = 18.199817
$ws.Cells.Item(10, 15).Value = 0.849784628791665
$ws.Cells.Item(10, 16).Value = 0.8497846287916652
$ws.Cells.Item(10, 17).Value = 8.351768622581
$ws.Cells.Item(10, 18).Value = 75.16591760322899
$ws.Cells.Item(10, 19).Value = 0.0004091790484062309
$ws.Cells.Item(10, 20).Value = 0.000409179048406231

$ws.Cells.Item(11, 7).Value = 1.376679
$ws.Cells.Item(11, 8).Value = 4.130037
$ws.Cells.Item(11, 9).Value = 0.0004815091195378001
$ws.Cells.Item(11, 10).Value = 0.0004815091195378002
$ws.Cells.Item(11, 15).Value = 0.1196497582104962
$ws.Cells.Item(11, 16).Value = 0.1196497582104962
$ws.Cells.Item(11, 17).Value = 1.175929832648
$ws.Cells.Item(11, 18).Value = 10.583368493832
$ws.Cells.Item(11, 19).Value = 0.00005761244972884668
$ws.Cells.Item(11, 20).Value = 0.0000576124497288467

$ws.Cells.Item(12, 7).Value = 1.376679
$ws.Cells.Item(12, 8).Value = 4.130037
$ws.Cells.Item(12, 9).Value = 0.0004815091195378001
$ws.Cells.Item(12, 10).Value = 0.0004815091195378002
$ws.Cells.Item(12, 13).Value = 0.1824346666666667
$ws.Cells.Item(12, 14).Value = 0.547304
$ws.Cells.Item(12, 15).Value = 0.02555468148257719
$ws.Cells.Item(12, 16).Value = 0.02555468148257719
$ws.Cells.Item(12, 17).Value = 0.251153974472
$ws.Cells.Item(12, 18).Value = 2.260385770248
$ws.Cells.Item(12, 19).Value = 0.00001230481218074466
$ws.Cells.Item(12, 20).Value = 0.00001230481218074467

$ws.Cells.Item(13, 7).Value = 1.376679
$ws.Cells.Item(13, 8).Value = 4.130037
$ws.Cells.Item(13, 9).Value = 0.0004815091195378001
$ws.Cells.Item(13, 10).Value = 0.0004815091195378002
$ws.Cells.Item(13, 13).Value = 0.035773
$ws.Cells.Item(13, 14).Value = 0.107319
$ws.Cells.Item(13, 15).Value = 0.005010931515261538
$ws.Cells.Item(13, 16).Value = 0.005010931515261539
$ws.Cells.Item(13, 17).Value = 0.049247937867
$ws.Cells.Item(13, 18).Value = 0.443231440803
$ws.Cells.Item(13, 19).Value = 0.000002412809221977798
$ws.Cells.Item(13, 20).Value = 0.000002412809221977798

$ws.Cells.Item(14, 7).Value = 164.8447596666666
$ws.Cells.Item(14, 8).Value = 494.534279
$ws.Cells.Item(14, 9).Value = 0.05765632735555414
$ws.Cells.Item(14, 10).Value = 0.05765632735555416
$ws.Cells.Item(14, 13).Value = 6.066605666666667
$ws.Cells.Item(14, 14).Value = 18.199817
$ws.Cells.Item(14, 15).Value = 0.849784628791665
$ws.Cells.Item(14, 16).Value = 0.8497846287916652
$ws.Cells.Item(14, 17).Value = 1000.048153114105
$ws.Cells.Item(14, 18).Value = 9000.433378026943
$ws.Cells.Item(14, 19).Value = 0.0489954607393303
$ws.Cells.Item(14, 20).Value = 0.04899546073933032

$ws.Cells.Item(15, 7).Value = 164.8447596666666
$ws.Cells.Item(15, 8).Value = 494.534279
$ws.Cells.Item(15, 9).Value = 0.05765632735555414
$ws.Cells.Item(15, 10).Value = 0.05765632735555416
$ws.Cells.Item(15, 15).Value = 0.1196497582104962
$ws.Cells.Item(15, 16).Value = 0.1196497582104962
$ws.Cells.Item(15, 17).Value = 140.8068770190604
$ws.Cells.Item(15, 18).Value = 1267.261893171544
$ws.Cells.Item(15, 19).Value = 0.006898565627397269
$ws.Cells.Item(15, 20).Value = 0.006898565627397273

$ws.Cells.Item(16, 7).Value = 164.8447596666666
$ws.Cells.Item(16, 8).Value = 494.534279
$ws.Cells.Item(16, 9).Value = 0.05765632735555414
$ws.Cells.Item(16, 10).Value = 0.05765632735555416
$ws.Cells.Item(16, 13).Value = 0.1824346666666667
$ws.Cells.Item(16, 14).Value = 0.547304
$ws.Cells.Item(16, 15).Value = 0.02555468148257719
$ws.Cells.Item(16, 16).Value = 0.02555468148257719
$ws.Cells.Item(16, 17).Value = 30.07339878153511
$ws.Cells.Item(16, 18).Value = 270.660589033816
$ws.Cells.Item(16, 19).Value = 0.001473389081026388
$ws.Cells.Item(16, 20).Value = 0.001473389081026388

$ws.Cells.Item(17, 7).Value = 164.8447596666666
$ws.Cells.Item(17, 8).Value = 494.534279
$ws.Cells.Item(17, 9).Value = 0.05765632735555414
$ws.Cells.Item(17, 10).Value = 0.05765632735555416
$ws.Cells.Item(17, 13).Value = 0.035773
$ws.Cells.Item(17, 14).Value = 0.107319
$ws.Cells.Item(17, 15).Value = 0.005010931515261538
$ws.Cells.Item(17, 16).Value = 0.005010931515261539
$ws.Cells.Item(17, 17).Value = 5.896991587555666
$ws.Cells.Item(17, 18).Value = 53.072924288001
$ws.Cells.Item(17, 19).Value = 0.0002889119078001822
$ws.Cells.Item(17, 20).Value = 0.0002889119078001823

$ws.Cells.Item(18, 5).Value = 3
$ws.Cells.Item(18, 6).Value = 1
$ws.Cells.Item(18, 7).Value = 1.431418
$ws.Cells.Item(18, 8).Value = 4.294254
$ws.Cells.Item(18, 9).Value = 0.0005006547066313635
$ws.Cells.Item(18, 10).Value = 0.0005006547066313636
$ws.Cells.Item(18, 13).Value = 6.066605666666667
$ws.Cells.Item(18, 14).Value = 18.199817
$ws.Cells.Item(18, 15).Value = 0.849784628791665
$ws.Cells.Item(18, 16).Value = 0.8497846287916652
$ws.Cells.Item(18, 17).Value = 8.683848550168667
$ws.Cells.Item(18, 18).Value = 78.154636951518
$ws.Cells.Item(18, 19).Value = 0.0004254486740275332
$ws.Cells.Item(18, 20).Value = 0.0004254486740275333

$ws.Cells.Item(19, 5).Value = 3
$ws.Cells.Item(19, 6).Value = 1
$ws.Cells.Item(19, 7).Value = 1.431418
$ws.Cells.Item(19, 8).Value = 4.294254
$ws.Cells.Item(19, 9).Value = 0.0005006547066313635
$ws.Cells.Item(19, 10).Value = 0.0005006547066313636
$ws.Cells.Item(19, 15).Value = 0.1196497582104962
$ws.Cells.Item(19, 16).Value = 0.1196497582104962
$ws.Cells.Item(19, 17).Value = 1.222686718682666
$ws.Cells.Item(19, 18).Value = 11.004180468144
$ws.Cells.Item(19, 19).Value = 0.00005990321459538953
$ws.Cells.Item(19, 20).Value = 0.00005990321459538956

$ws.Cells.Item(20, 5).Value = 3
$ws.Cells.Item(20, 6).Value = 1
$ws.Cells.Item(20, 7).Value = 1.431418
$ws.Cells.Item(20, 8).Value = 4.294254
$ws.Cells.Item(20, 9).Value = 0.0005006547066313635
$ws.Cells.Item(20, 10).Value = 0.0005006547066313636
$ws.Cells.Item(20, 13).Value = 0.1824346666666667
$ws.Cells.Item(20, 14).Value = 0.547304
$ws.Cells.Item(20, 15).Value = 0.02555468148257719
$ws.Cells.Item(20, 16).Value = 0.02555468148257719
$ws.Cells.Item(20, 17).Value = 0.2611402656906667
$ws.Cells.Item(20, 18).Value = 2.350262391216
$ws.Cells.Item(20, 19).Value = 0.00001279407156071762
$ws.Cells.Item(20, 20).Value = 0.00001279407156071762

$ws.Cells.Item(21, 5).Value = 3
$ws.Cells.Item(21, 6).Value = 1
$ws.Cells.Item(21, 7).Value = 1.431418
$ws.Cells.Item(21, 8).Value = 4.294254
$ws.Cells.Item(21, 9).Value = 0.0005006547066313635
$ws.Cells.Item(21, 10).Value = 0.0005006547066313636
$ws.Cells.Item(21, 13).Value = 0.035773
$ws.Cells.Item(21, 14).Value = 0.107319
$ws.Cells.Item(21, 15).Value = 0.005010931515261538
$ws.Cells.Item(21, 16).Value = 0.005010931515261539
$ws.Cells.Item(21, 17).Value = 0.051206116114
$ws.Cells.Item(21, 18).Value = 0.4608550450260001
$ws.Cells.Item(21, 19).Value = 0.000002412809221977798
$ws.Cells.Item(21, 20).Value = 0.000002412809221977798

